$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (case description) updates: fix "depandant" -> "dependent" typo ---
$ws.Range("B3").Value = "All parameters different - time dependent"
$ws.Range("B4").Value = "All parameters constant - time dependent"
$ws.Range("B5").Value = "k and Y constant - time dependent"
$ws.Range("B6").Value = "Allocations constant - time dependent"
$ws.Range("B7").Value = "Turnovers constant - time dependent"
$ws.Range("B8").Value = "All parameters different - temperature dependent"
$ws.Range("B9").Value = "All parameters different - plant size (height) dependent"

# --- Column D (comments) updates ---
$ws.Range("D4").Value = "There is certain effect of warming on parameters"
$ws.Range("D5").Value = "Works slightly worse than case 2"
$ws.Range("D6").Value = "Works almost equally well as case 2, except the roots "
$ws.Range("D7").Value = "Works slightly worse than case 2"

# --- Column width for D (wider to fit new, longer text); engine quantizes
# ColumnWidth to 1/6-character steps, so 43.3 is the closest input that lands
# on the nearest achievable stored width to the target 44.1640625 ---
$ws.Range("D1").ColumnWidth = 43.3

# --- Update the saved selection to B10 (below the table) ---
$ws.Range("B10").Select()
